$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "26.896.75"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +1.01%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.843.77"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.18%  "

$ws.Range("E4").Value = "  -0.41%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "308.91"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.96%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.005"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.36%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4762"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.72%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3679"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +2.33%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07199"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.86%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.9260"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +2.55%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "19.62"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.96%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07605"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -2.59%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.894.37"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +3.43%  "

$ws.Range("E14").Value = "  +1.00%  "

$ws.Range("E15").Value = "  +1.02%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "88.53"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.37%  "

$ws.Range("E17").Value = "  -0.17%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008655"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.04%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "1.005"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.32%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "26.944.03"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.99%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "14.54"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.43%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.033"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.36%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.65"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.79%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "1.919"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.01%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "152.06"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.10%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "18.14"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.25%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.999"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +1.14%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "114.26"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.57%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "4.945"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.74%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.08853"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.53%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.305"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +5.11%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.7486"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +2.34%  "

$ws.Range("E33").Value = "  +4.06%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "2.766"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.46%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "4.482"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.92%  "

$ws.Range("E36").Value = "  +0.93%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.05259"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.83%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.01947"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.85%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.958"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.30%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.5221"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +3.05%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "6.917"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.24%  "

$ws.Range("E42").Value = "  +1.01%  "

$ws.Range("E43").Value = "  +2.80%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "10.51"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +4.82%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.4696"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.26%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.005"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.37%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "101.82"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.75%  "

$ws.Range("E48").Value = "  +2.95%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "65.36"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +2.39%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.06029"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.8849"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +3.94%  "
